# Scheduled runner update: refresh computed profit-margin figures on several
# sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR) to reflect the latest market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 500
$ws.Range("I31").Value = 500
$ws.Range("K31").Value = 1500
$ws.Range("M31").Value = -1270

$ws.Range("H43").Value = 3125.75
$ws.Range("I43").Value = 3167
$ws.Range("K43").Value = 3167
$ws.Range("M43").Value = -3098

$ws.Range("H62").Value = 3500
$ws.Range("J62").Value = 3500
$ws.Range("L62").Value = 3500
$ws.Range("N62").Value = -4748

$ws.Range("H65").Value = 3500
$ws.Range("J65").Value = 3500
$ws.Range("L65").Value = 17500
$ws.Range("N65").Value = -23740

$ws.Range("H111").Value = 2517.4167
$ws.Range("I111").Value = 2517.4167
$ws.Range("K111").Value = 7552.250100000001
$ws.Range("M111").Value = -4485.250100000001

$ws.Range("H138").Value = 1318.871
$ws.Range("I138").Value = 1383.2667
$ws.Range("J138").Value = 1258.5
$ws.Range("K138").Value = 4149.800099999999
$ws.Range("L138").Value = 3775.5
$ws.Range("M138").Value = 990.1999000000005
$ws.Range("N138").Value = -14055.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 22000
$ws.Range("I37").Value = 4000
$ws.Range("K37").Value = 4000
$ws.Range("M37").Value = -3727

$ws.Range("H74").Value = 3959.6667
$ws.Range("I74").Value = 3069.625
$ws.Range("J74").Value = 4976.857
$ws.Range("K74").Value = 3069.625
$ws.Range("L74").Value = 4976.857
$ws.Range("M74").Value = -2195.625
$ws.Range("N74").Value = -6724.857

$ws.Range("H77").Value = 3959.6667
$ws.Range("I77").Value = 3069.625
$ws.Range("J77").Value = 4976.857
$ws.Range("K77").Value = 15348.125
$ws.Range("L77").Value = 24884.285
$ws.Range("M77").Value = -10980.125
$ws.Range("N77").Value = -33620.285

$ws.Range("H80").Value = 69396.664
$ws.Range("J80").Value = 69396.664
$ws.Range("L80").Value = 69396.664
$ws.Range("N80").Value = -71392.664

$ws.Range("H83").Value = 69396.664
$ws.Range("J83").Value = 69396.664
$ws.Range("L83").Value = 208189.992
$ws.Range("N83").Value = -218173.992

$ws.Range("H110").Value = 8254.4375
$ws.Range("I110").Value = 6889.5713
$ws.Range("K110").Value = 6889.5713
$ws.Range("M110").Value = -4844.5713

$ws.Range("H139").Value = 99999.5
$ws.Range("J139").Value = 99999.5
$ws.Range("L139").Value = 99999.5
$ws.Range("N139").Value = -110279.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 46037
$ws.Range("J35").Value = 60074
$ws.Range("L35").Value = 60074
$ws.Range("N35").Value = -60694

$ws.Range("H59").Value = 130000
$ws.Range("I59").Value = 130000
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 130000
$ws.Range("L59").Value = 0
$ws.Range("M59").Value = -129153
$ws.Range("N59").ClearContents()

$ws.Range("H86").Value = 5032.4736
$ws.Range("I86").Value = 2622
$ws.Range("K86").Value = 2622
$ws.Range("M86").Value = -1499

$ws.Range("H89").Value = 5032.4736
$ws.Range("I89").Value = 2622
$ws.Range("K89").Value = 13110
$ws.Range("M89").Value = -7494

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2987.1428
$ws.Range("I31").Value = 2984.6
$ws.Range("K31").Value = 2984.6
$ws.Range("M31").Value = -2689.6

$ws.Range("H34").Value = 2987.1428
$ws.Range("I34").Value = 2984.6
$ws.Range("K34").Value = 2984.6
$ws.Range("M34").Value = -2782.6

$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("N45").ClearContents()

$ws.Range("H51").Value = 44999.5
$ws.Range("J51").Value = 50000
$ws.Range("L51").Value = 50000
$ws.Range("N51").Value = -51472

$ws.Range("H58").Value = 2073.2307
$ws.Range("I58").Value = 1515.3334
$ws.Range("J58").Value = 3328.5
$ws.Range("K58").Value = 1515.3334
$ws.Range("L58").Value = 3328.5
$ws.Range("M58").Value = -1312.3334
$ws.Range("N58").Value = -3734.5

$ws.Range("H59").Value = 104
$ws.Range("I59").Value = 104
$ws.Range("K59").Value = 104
$ws.Range("M59").Value = 1041

$ws.Range("H60").Value = 96.5
$ws.Range("I60").Value = 96.5
$ws.Range("K60").Value = 96.5
$ws.Range("M60").Value = 414.5

$ws.Range("H61").Value = 44999.5
$ws.Range("J61").Value = 50000
$ws.Range("L61").Value = 50000
$ws.Range("N61").Value = -50696

$ws.Range("H68").Value = 70000
$ws.Range("J68").Value = 70000
$ws.Range("L68").Value = 70000
$ws.Range("N68").Value = -71498

$ws.Range("H71").Value = 70000
$ws.Range("J71").Value = 70000
$ws.Range("L71").Value = 210000
$ws.Range("N71").Value = -217488

$ws.Range("H136").Value = 2073.2307
$ws.Range("I136").Value = 1515.3334
$ws.Range("J136").Value = 3328.5
$ws.Range("K136").Value = 4546.0002
$ws.Range("L136").Value = 9985.5
$ws.Range("M136").Value = -1996.0002
$ws.Range("N136").Value = -15085.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 10794.333
$ws.Range("I43").Value = 837.1111
$ws.Range("J43").Value = 40666
$ws.Range("K43").Value = 837.1111
$ws.Range("L43").Value = 40666
$ws.Range("M43").Value = -686.1111
$ws.Range("N43").Value = -40968

$ws.Range("H46").Value = 26188.611
$ws.Range("I46").Value = 6166.3335
$ws.Range("J46").Value = 30193.066
$ws.Range("K46").Value = 6166.3335
$ws.Range("L46").Value = 30193.066
$ws.Range("M46").Value = -6010.3335
$ws.Range("N46").Value = -30505.066

$ws.Range("H138").Value = 83000
$ws.Range("J138").Value = 83000
$ws.Range("L138").Value = 83000
$ws.Range("N138").Value = -93280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3898.6667
$ws.Range("I22").Value = 899
$ws.Range("J22").Value = 5398.5
$ws.Range("K22").Value = 899
$ws.Range("L22").Value = 5398.5
$ws.Range("M22").Value = -604
$ws.Range("N22").Value = -5988.5

$ws.Range("H27").Value = 3898.6667
$ws.Range("I27").Value = 899
$ws.Range("J27").Value = 5398.5
$ws.Range("K27").Value = 899
$ws.Range("L27").Value = 5398.5
$ws.Range("M27").Value = -792
$ws.Range("N27").Value = -5612.5

$ws.Range("H46").Value = 1374.75
$ws.Range("I46").Value = 1102.9412
$ws.Range("J46").Value = 2915
$ws.Range("K46").Value = 1102.9412
$ws.Range("L46").Value = 2915
$ws.Range("M46").Value = -914.9412
$ws.Range("N46").Value = -3291

$ws.Range("H93").Value = 3040.5454
$ws.Range("I93").Value = 2855.7144
$ws.Range("K93").Value = 2855.7144
$ws.Range("M93").Value = -1607.7144

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()

$ws.Range("H54").Value = 21352.092
$ws.Range("I54").Value = 9500
$ws.Range("J54").Value = 28124.715
$ws.Range("K54").Value = 9500
$ws.Range("L54").Value = 28124.715
$ws.Range("M54").Value = -8980
$ws.Range("N54").Value = -29164.715
